# Auto update: 2025-12-05 19:04:26
# Update final score (K) and MACRO_SCORE (N) values for rows 2 and 3
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 56.7
$ws.Range("N2").Value = 51.15965480231979

$ws.Range("K3").Value = 53.3
$ws.Range("N3").Value = 51.15965480231979
